$wb = $excel.ActiveWorkbook

$wsProc = $wb.Worksheets.Item("SetsEditor- Proc")
$wsComm = $wb.Worksheets.Item("VEDA_Sets-Comm")

# --- VEDA_Sets-Comm (sheet2): tidy up B9 formula and rename D9 string ---
$wsComm.Activate()
$wsComm.Range("B9").Formula = "=B10&"",""&B11&"",""&B12&"",""&B13&"",""&B15"
$wsComm.Range("D9").Value = "FinalEnergy"
$wsComm.Range("B10").Select()

# --- SetsEditor- Proc (sheet1): insert a new row for IND_ElecGenDist ---
$wsProc.Activate()
$wsProc.Rows.Item(24).Insert()
$wsProc.Range("E24").Value = "IND_ElecGenDist"
$wsProc.Range("F24").Value = "IND_ElecGenDist"
$wsProc.Range("H24").Value = "E*"
$wsProc.Range("J24").Value = "I??ELC"
$wsProc.Range("M24").Value = "L2"
$wsProc.Range("N24").Value = "SubSector"
$wsProc.Range("J24").Select()

$wsProc.Activate()
